{"js": "// Replace the 25 division-problem texts in the practice table with the\n// new set of problems, preserving all run formatting (font/size) since\n// we only change the text content of each matched run via search+replace.\nconst replacements = [\n  [\"572\u00f72=\", \"375\u00f78=\"],\n  [\"890\u00f76=\", \"621\u00f72=\"],\n  [\"479\u00f75=\", \"973\u00f73=\"],\n  [\"758\u00f72=\", \"311\u00f79=\"],\n  [\"261\u00f75=\", \"939\u00f75=\"],\n  [\"577\u00f74=\", \"573\u00f75=\"],\n  [\"903\u00f79=\", \"564\u00f78=\"],\n  [\"545\u00f74=\", \"797\u00f75=\"],\n  [\"576\u00f78=\", \"954\u00f73=\"],\n  [\"744\u00f72=\", \"230\u00f76=\"],\n  [\"793\u00f76=\", \"148\u00f76=\"],\n  [\"528\u00f78=\", \"290\u00f75=\"],\n  [\"108\u00f79=\", \"340\u00f74=\"],\n  [\"749\u00f76=\", \"913\u00f76=\"],\n  [\"224\u00f74=\", \"658\u00f78=\"],\n  [\"260\u00f78=\", \"925\u00f75=\"],\n  [\"607\u00f73=\", \"286\u00f76=\"],\n  [\"927\u00f78=\", \"267\u00f72=\"],\n  [\"923\u00f78=\", \"902\u00f79=\"],\n  [\"350\u00f72=\", \"380\u00f78=\"],\n  [\"121\u00f77=\", \"807\u00f79=\"],\n  [\"219\u00f73=\", \"294\u00f73=\"],\n  [\"991\u00f75=\", \"706\u00f72=\"],\n  [\"683\u00f75=\", \"980\u00f73=\"],\n  [\"254\u00f74=\", \"732\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem texts in the practice table with the\n# new set of problems. Using Find/Replace on $d.Content preserves the\n# existing run formatting (font/size) of each matched text run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"572\u00f72=\", \"375\u00f78=\"),\n  @(\"890\u00f76=\", \"621\u00f72=\"),\n  @(\"479\u00f75=\", \"973\u00f73=\"),\n  @(\"758\u00f72=\", \"311\u00f79=\"),\n  @(\"261\u00f75=\", \"939\u00f75=\"),\n  @(\"577\u00f74=\", \"573\u00f75=\"),\n  @(\"903\u00f79=\", \"564\u00f78=\"),\n  @(\"545\u00f74=\", \"797\u00f75=\"),\n  @(\"576\u00f78=\", \"954\u00f73=\"),\n  @(\"744\u00f72=\", \"230\u00f76=\"),\n  @(\"793\u00f76=\", \"148\u00f76=\"),\n  @(\"528\u00f78=\", \"290\u00f75=\"),\n  @(\"108\u00f79=\", \"340\u00f74=\"),\n  @(\"749\u00f76=\", \"913\u00f76=\"),\n  @(\"224\u00f74=\", \"658\u00f78=\"),\n  @(\"260\u00f78=\", \"925\u00f75=\"),\n  @(\"607\u00f73=\", \"286\u00f76=\"),\n  @(\"927\u00f78=\", \"267\u00f72=\"),\n  @(\"923\u00f78=\", \"902\u00f79=\"),\n  @(\"350\u00f72=\", \"380\u00f78=\"),\n  @(\"121\u00f77=\", \"807\u00f79=\"),\n  @(\"219\u00f73=\", \"294\u00f73=\"),\n  @(\"991\u00f75=\", \"706\u00f72=\"),\n  @(\"683\u00f75=\", \"980\u00f73=\"),\n  @(\"254\u00f74=\", \"732\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
